$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new server data row (row 2)
# Columns: A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
$ws.Range("C2").NumberFormat = "@"

# Set values in an order that reproduces the shared-string table ordering
$ws.Range("B2").Value = "000106001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "MasterServer_1"
$ws.Range("C2").Value = "MasterServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# Update the active selection to H3
$ws.Range("H3").Select()
